$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.245.61'
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").Value = '1.860.66'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7028'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.65'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08140'
$ws.Range("E8").Value = '  +9.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3028'
$ws.Range("E9").Value = '  -0.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.20'
$ws.Range("E10").Value = '  -0.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08162'
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").Value = '1.846.24'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.160'
$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7053'
$ws.Range("E14").Value = '  -2.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.04'
$ws.Range("E15").Value = '  +0.38%  '

$ws.Range("D16").Value = '29.260.69'
$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.772'
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007837'
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.31'
$ws.Range("E19").Value = '  +1.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.52'
$ws.Range("E20").Value = '  -0.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'

$ws.Range("D22").Value = '2.111.82'
$ws.Range("E22").Value = '  +1.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.411'
$ws.Range("E24").Value = '  -1.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.23'
$ws.Range("E25").Value = '  +0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.953'
$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.06'
$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.966'
$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.431'
$ws.Range("E30").Value = '  +2.52%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.480'
$ws.Range("E31").Value = '  -0.72%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.394'
$ws.Range("E32").Value = '  -2.47%  '

$ws.Range("E33").Value = '  +2.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05179'
$ws.Range("E34").Value = '  +0.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.167'
$ws.Range("E35").Value = '  -1.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7058'
$ws.Range("E36").Value = '  +1.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9982'
$ws.Range("E37").Value = '  -1.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.673'
$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01839'
$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.729'
$ws.Range("E40").Value = '  +1.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9207'
$ws.Range("E41").Value = '  -2.22%  '

$ws.Range("D42").Value = '1.135.44'
$ws.Range("E42").Value = '  +5.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4265'
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.12'
$ws.Range("E45").Value = '  +0.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.20'
$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.763'
$ws.Range("E48").Value = '  +1.45%  '

$ws.Range("D49").Value = '2.007.58'
$ws.Range("E49").Value = '  +0.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.155'
$ws.Range("E50").Value = '  +0.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.942'
$ws.Range("E51").Value = '  -1.13%  '
